# ev branch under test
# Update supply_temperature (column Q) values and adjust the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 52
$ws.Range("Q3").Value = 52
$ws.Range("Q4").Value = 52

$ws.Range("Q5").Select()
